# Reorder the player roster rows (A5:C16) to match the updated sheet.
# Rows 2-4 and 17-19 stay unchanged; rows 5-16 are re-sequenced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Dillon Brooks", "SG,SF,PF", "Houston Rockets"),
    @("Paolo Banchero", "SF,PF", "Orlando Magic"),
    @("Ayo Dosunmu", "PG,SG,SF", "Chicago Bulls"),
    @("Pascal Siakam", "SF,PF,C", "Indiana Pacers"),
    @("Rudy Gobert", "C", "Minnesota Timberwolves"),
    @("Chet Holmgren", "PF,C", "Oklahoma City Thunder"),
    @("Jalen Green", "PG,SG", "Houston Rockets"),
    @("Aaron Wiggins", "SG,SF", "Oklahoma City Thunder"),
    @("Brandin Podziemski", "SG", "Golden State Warriors"),
    @("Deni Avdija", "SF,PF", "Portland Trail Blazers"),
    @("Nikola Jokic", "C", "Denver Nuggets"),
    @("Jaylen Brown", "SG,SF", "Boston Celtics")
)

$startRow = 5
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
